# WBS.xlsx - dcm and pgm commit on 20230918
# Update the "release" (リリース) milestone row: set both the planned
# start date (予定開始日, column D) and planned end date (予定完了日,
# column E) to 2023-09-25.
#
# The cells already carry a date-formatted style (s="2", numFmtId 56),
# so we write the underlying date serial value directly rather than a
# date string, to avoid Excel minting a brand-new number-format style
# for the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D21").Value = 45194
$ws.Range("E21").Value = 45194
